$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Embedded target data for rows 2-193 (numeric columns A,B,C,D)
# Format per line: row,A,B,C,D
$rowData = @"
2,45969.0,333,567,900
3,45969.01041666666,329,554,883
4,45969.02083333334,328,573,901
5,45969.03125,326,0,326
6,45969.04166666666,328,588,916
7,45969.05208333334,329,589,918
8,45969.0625,330,590,920
9,45969.07291666666,328,588,916
10,45969.08333333334,0,594,594
11,45969.09375,327,592,919
12,45969.10416666666,328,593,921
13,45969.11458333334,326,0,326
14,45969.125,327,589,916
15,45969.13541666666,0,0,0
16,45969.14583333334,326,0,326
17,45969.15625,327,590,917
18,45969.16666666666,326,598,924
19,45969.17708333334,325,0,325
20,45969.1875,321,0,321
21,45969.19791666666,311,599,910
22,45969.20833333334,336,670,1006
23,45969.21875,340,674,1014
24,45969.22916666666,0,0,0
25,45969.23958333334,0,679,679
26,45969.25,483,714,1197
27,45969.26041666666,433,711,1144
28,45969.27083333334,0,650,650
29,45969.28125,436,648,1084
30,45969.29166666666,580,807,1387
31,45969.30208333334,518,772,1290
32,45969.3125,509,765,1274
33,45969.32291666666,512,761,1273
34,45969.33333333334,375,703,1078
35,45969.34375,374,700,1074
36,45969.35416666666,447,791,1238
37,45969.36458333334,449,806,1255
38,45969.375,293,763,1056
39,45969.38541666666,279,771,1050
40,45969.39583333334,276,778,1054
41,45969.40625,274,780,1054
42,45969.41666666666,299,709,1008
43,45969.42708333334,308,704,1012
44,45969.4375,306,694,1000
45,45969.44791666666,309,705,1014
46,45969.45833333334,216,722,938
47,45969.46875,208,724,932
48,45969.47916666666,209,725,934
49,45969.48958333334,207,0,207
50,45969.5,227,604,831
51,45969.51041666666,231,597,828
52,45969.52083333334,234,0,234
53,45969.53125,215,591,806
54,45969.54166666666,319,610,929
55,45969.55208333334,298,609,907
56,45969.5625,294,0,294
57,45969.57291666666,336,0,336
58,45969.58333333334,494,799,1293
59,45969.59375,512,806,1318
60,45969.60416666666,526,808,1334
61,45969.61458333334,536,0,536
62,45969.625,617,990,1607
63,45969.63541666666,603,1013,1616
64,45969.64583333334,602,0,602
65,45969.65625,0,984,984
66,45969.66666666666,637,1008,1645
67,45969.67708333334,648,1037,1685
68,45969.6875,664,0,664
69,45969.69791666666,663,895,1558
70,45969.70833333334,686,1115,1801
71,45969.71875,687,1123,1810
72,45969.72916666666,684,1107,1791
73,45969.73958333334,683,1054,1737
74,45969.75,650,1089,1739
75,45969.76041666666,767,1079,1846
76,45969.77083333334,762,1065,1827
77,45969.78125,769,965,1734
78,45969.79166666666,648,1034,1682
79,45969.80208333334,0,1037,1037
80,45969.8125,731,1054,1785
81,45969.82291666666,0,1056,1056
82,45969.83333333334,682,1044,1726
83,45969.84375,670,1042,1712
84,45969.85416666666,671,0,671
85,45969.86458333334,669,1043,1712
86,45969.875,599,1046,1645
87,45969.88541666666,594,0,594
88,45969.89583333334,609,1050,1659
89,45969.90625,598,1091,1689
90,45969.91666666666,358,765,1123
91,45969.92708333334,349,742,1091
92,45969.9375,351,0,351
93,45969.94791666666,356,864,1220
94,45969.95833333334,380,844,1224
95,45969.96875,379,840,1219
96,45969.97916666666,377,841,1218
97,45969.98958333334,381,842,1223
98,45970.0,378,627,1005
99,45970.01041666666,377,614,991
100,45970.02083333334,0,0,0
101,45970.03125,368,613,981
102,45970.04166666666,358,631,989
103,45970.05208333334,357,632,989
104,45970.0625,0,0,0
105,45970.07291666666,359,0,359
106,45970.08333333334,362,635,997
107,45970.09375,359,591,950
108,45970.10416666666,360,589,949
109,45970.11458333334,359,590,949
110,45970.125,0,630,630
111,45970.13541666666,361,632,993
112,45970.14583333334,0,0,0
113,45970.15625,360,0,360
114,45970.16666666666,363,650,1013
115,45970.17708333334,361,0,361
116,45970.1875,358,0,358
117,45970.19791666666,360,651,1011
118,45970.20833333334,361,729,1090
119,45970.21875,360,740,1100
120,45970.22916666666,359,756,1115
121,45970.23958333334,365,757,1122
122,45970.25,396,693,1089
123,45970.26041666666,395,700,1095
124,45970.27083333334,396,0,396
125,45970.28125,397,0,397
126,45970.29166666666,407,674,1081
127,45970.30208333334,0,0,0
128,45970.3125,0,0,0
129,45970.32291666666,0,0,0
130,45970.33333333334,0,0,0
131,45970.34375,0,0,0
132,45970.35416666666,0,0,0
133,45970.36458333334,0,0,0
134,45970.375,0,0,0
135,45970.38541666666,0,0,0
136,45970.39583333334,0,0,0
137,45970.40625,0,0,0
138,45970.41666666666,0,0,0
139,45970.42708333334,0,0,0
140,45970.4375,0,0,0
141,45970.44791666666,0,0,0
142,45970.45833333334,0,0,0
143,45970.46875,0,0,0
144,45970.47916666666,0,0,0
145,45970.48958333334,0,0,0
146,45970.5,0,0,0
147,45970.51041666666,0,0,0
148,45970.52083333334,0,0,0
149,45970.53125,0,0,0
150,45970.54166666666,0,0,0
151,45970.55208333334,0,0,0
152,45970.5625,0,0,0
153,45970.57291666666,0,0,0
154,45970.58333333334,0,0,0
155,45970.59375,0,0,0
156,45970.60416666666,0,0,0
157,45970.61458333334,0,0,0
158,45970.625,0,0,0
159,45970.63541666666,0,0,0
160,45970.64583333334,0,0,0
161,45970.65625,0,0,0
162,45970.66666666666,0,0,0
163,45970.67708333334,0,0,0
164,45970.6875,0,0,0
165,45970.69791666666,0,0,0
166,45970.70833333334,0,0,0
167,45970.71875,0,0,0
168,45970.72916666666,0,0,0
169,45970.73958333334,0,0,0
170,45970.75,0,0,0
171,45970.76041666666,0,0,0
172,45970.77083333334,0,0,0
173,45970.78125,0,0,0
174,45970.79166666666,0,0,0
175,45970.80208333334,0,0,0
176,45970.8125,0,0,0
177,45970.82291666666,0,0,0
178,45970.83333333334,0,0,0
179,45970.84375,0,0,0
180,45970.85416666666,0,0,0
181,45970.86458333334,0,0,0
182,45970.875,0,0,0
183,45970.88541666666,0,0,0
184,45970.89583333334,0,0,0
185,45970.90625,0,0,0
186,45970.91666666666,0,0,0
187,45970.92708333334,0,0,0
188,45970.9375,0,0,0
189,45970.94791666666,0,0,0
190,45970.95833333334,0,0,0
191,45970.96875,0,0,0
192,45970.97916666666,0,0,0
193,45970.98958333334,0,0,0
"@

# Embedded target data for rows 2-193 (text column F - "Lookup")
# Format per line: row,text
$fData = @"
2,08.11.20251
3,08.11.20252
4,08.11.20253
5,08.11.20254
6,08.11.20255
7,08.11.20256
8,08.11.20257
9,08.11.20258
10,08.11.20259
11,08.11.202510
12,08.11.202511
13,08.11.202512
14,08.11.202513
15,08.11.202514
16,08.11.202515
17,08.11.202516
18,08.11.202517
19,08.11.202518
20,08.11.202519
21,08.11.202520
22,08.11.202521
23,08.11.202522
24,08.11.202523
25,08.11.202524
26,08.11.202525
27,08.11.202526
28,08.11.202527
29,08.11.202528
30,08.11.202529
31,08.11.202530
32,08.11.202531
33,08.11.202532
34,08.11.202533
35,08.11.202534
36,08.11.202535
37,08.11.202536
38,08.11.202537
39,08.11.202538
40,08.11.202539
41,08.11.202540
42,08.11.202541
43,08.11.202542
44,08.11.202543
45,08.11.202544
46,08.11.202545
47,08.11.202546
48,08.11.202547
49,08.11.202548
50,08.11.202549
51,08.11.202550
52,08.11.202551
53,08.11.202552
54,08.11.202553
55,08.11.202554
56,08.11.202555
57,08.11.202556
58,08.11.202557
59,08.11.202558
60,08.11.202559
61,08.11.202560
62,08.11.202561
63,08.11.202562
64,08.11.202563
65,08.11.202564
66,08.11.202565
67,08.11.202566
68,08.11.202567
69,08.11.202568
70,08.11.202569
71,08.11.202570
72,08.11.202571
73,08.11.202572
74,08.11.202573
75,08.11.202574
76,08.11.202575
77,08.11.202576
78,08.11.202577
79,08.11.202578
80,08.11.202579
81,08.11.202580
82,08.11.202581
83,08.11.202582
84,08.11.202583
85,08.11.202584
86,08.11.202585
87,08.11.202586
88,08.11.202587
89,08.11.202588
90,08.11.202589
91,08.11.202590
92,08.11.202591
93,08.11.202592
94,08.11.202593
95,08.11.202594
96,08.11.202595
97,08.11.202596
98,09.11.20251
99,09.11.20252
100,09.11.20253
101,09.11.20254
102,09.11.20255
103,09.11.20256
104,09.11.20257
105,09.11.20258
106,09.11.20259
107,09.11.202510
108,09.11.202511
109,09.11.202512
110,09.11.202513
111,09.11.202514
112,09.11.202515
113,09.11.202516
114,09.11.202517
115,09.11.202518
116,09.11.202519
117,09.11.202520
118,09.11.202521
119,09.11.202522
120,09.11.202523
121,09.11.202524
122,09.11.202525
123,09.11.202526
124,09.11.202527
125,09.11.202528
126,09.11.202529
127,09.11.202530
128,09.11.202531
129,09.11.202532
130,09.11.202533
131,09.11.202534
132,09.11.202535
133,09.11.202536
134,09.11.202537
135,09.11.202538
136,09.11.202539
137,09.11.202540
138,09.11.202541
139,09.11.202542
140,09.11.202543
141,09.11.202544
142,09.11.202545
143,09.11.202546
144,09.11.202547
145,09.11.202548
146,09.11.202549
147,09.11.202550
148,09.11.202551
149,09.11.202552
150,09.11.202553
151,09.11.202554
152,09.11.202555
153,09.11.202556
154,09.11.202557
155,09.11.202558
156,09.11.202559
157,09.11.202560
158,09.11.202561
159,09.11.202562
160,09.11.202563
161,09.11.202564
162,09.11.202565
163,09.11.202566
164,09.11.202567
165,09.11.202568
166,09.11.202569
167,09.11.202570
168,09.11.202571
169,09.11.202572
170,09.11.202573
171,09.11.202574
172,09.11.202575
173,09.11.202576
174,09.11.202577
175,09.11.202578
176,09.11.202579
177,09.11.202580
178,09.11.202581
179,09.11.202582
180,09.11.202583
181,09.11.202584
182,09.11.202585
183,09.11.202586
184,09.11.202587
185,09.11.202588
186,09.11.202589
187,09.11.202590
188,09.11.202591
189,09.11.202592
190,09.11.202593
191,09.11.202594
192,09.11.202595
193,09.11.202596
"@

$rowLines = $rowData -split "`n" | Where-Object { $_.Trim() -ne "" }
$fLines = $fData -split "`n" | Where-Object { $_.Trim() -ne "" }

$n = $rowLines.Count
$arr = New-Object 'object[,]' $n,4
for ($i = 0; $i -lt $n; $i++) {
    $parts = $rowLines[$i].Split(",")
    $arr[$i,0] = [double]$parts[1]
    $arr[$i,1] = [double]$parts[2]
    $arr[$i,2] = [double]$parts[3]
    $arr[$i,3] = [double]$parts[4]
}
$ws.Range("A2:D193").Value = $arr

$farr = New-Object 'object[,]' $n,1
for ($i = 0; $i -lt $n; $i++) {
    $parts = $fLines[$i].Split(",")
    $farr[$i,0] = $parts[1]
}
$ws.Range("F2:F193").Value = $farr
